# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to reflect the latest generated output.
#
# Sheet "展览"   : row 2 (F2) 3357 -> 3359, row 5 (F5) 1475 -> 1493, row 6 (F6) 37 -> 41
# Sheet "全部类型": row 2 (F2) 3357 -> 3359, row 5 (F5) 1475 -> 1493, row 6 (F6) 37 -> 41

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 3359
    $ws.Range("F5").Value = 1493
    $ws.Range("F6").Value = 41
}
